# Third commit - after issue resolved
#
# Changes applied:
#   1. Cell B2 on the active sheet changes from "No" to "Yes"
#      (shared-string table collapses from 13 -> 12 unique strings because
#      "No" is no longer referenced anywhere).
#   2. The worksheet selection moves from C7 to the whole of row 4
#      (A4:XFD1048576), with A4 as the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the cell value.
$ws.Range("B2").Value = "Yes"

# 2. Select entire row 4 (A4:XFD1048576), active cell A4.
$ws.Range("A4:XFD1048576").Select() | Out-Null
